$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in newly-added "Price?" cells (column C) that previously had no value
$ws.Range("C5").Value = "1 royal per bar"
$ws.Range("C6").Value = "1.5 royals per bar"
$ws.Range("C22").Value = "8 coins per bag"
$ws.Range("C23").Value = "1.5 royals per bar"
$ws.Range("C35").Value = "2 royals per bar"
$ws.Range("C38").Value = "3 coins per sack"
$ws.Range("C39").Value = "8 coins per bar"
$ws.Range("C45").Value = "3 coins per sack"
$ws.Range("C46").Value = "7 coins per bar"

# Update the MEAT price text (chip/coin amounts increased)
$ws.Range("C26").Value = "9 chips per pound of pork. 9 chips per pound of beef. 12 chips per pound of lamb. 15 chips per pound of poultry. 1.2 coins per pound of veal. 3 coins per pound of venison."

# Update a few tooltip texts in column D
$ws.Range("D4").Value = "A barrel filled to the brim with beer. No minors allowed!"
$ws.Range("D8").Value = "A bushel of corn. Looks pretty tasty…"
$ws.Range("D16").Value = "A bundle of flax. Needs to be processed somehow."

# Restore the view/selection state left by the author after editing
try {
    $ws.Application.ActiveWindow.ScrollRow = 21
} catch {}
try {
    $ws.Range("C26").Select()
} catch {}
